$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old content entirely (rows 1, 3-10, cols A-G) and rebuild
# as a contiguous block A1:F9.
$ws.Cells.Clear()

# Data (regno, name, mobile, department, year, gender) per row, rows 2-9
$data = @(
    @(2024175052, "Dhanu",  9865234175, "IST", "3rd", "Female"),
    @(2024175053, "Trisha", 6589742351, "IST", "1st", "Male"),
    @(2024175054, "Mani",   9658413457, "CS",  "3rd", "Female"),
    @(2024175055, "Kavi",   7569841235, "CS",  "1nd", "Male"),
    @(2024175056, "Raja",   8654123975, "EEE", "2nd", "Male"),
    @(2024175057, "Priya",  9574621574, "EEE", "2nd", "Female"),
    @(2024175058, "Aadhi",  7896542856, "ECE", "4th", "Female"),
    @(2024175059, "Nila",   9865742563, "ECE", "4th", "Female")
)

# Fill column E (year) first: header then each row top-to-bottom.
$ws.Cells.Item(1, 5).Value = "year"
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Then the A1 and C1 headers.
$ws.Cells.Item(1, 1).Value = "regno"
$ws.Cells.Item(1, 3).Value = "mobile"

# Remaining header cells.
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 4).Value = "department"
$ws.Cells.Item(1, 6).Value = "gender"

# Remaining data cells.
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

$ws.Range("D12").Select()
